# Add data for 2022-05-04
# - Update "through" date label from April 25 to April 26 (sheet name + header cell)
# - Update carjacking counts for the newly-included day

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet and update the column header text to reflect the new cutoff date
$ws.Name = "Through 2022-04-26"
$ws.Range("B1").Value = "April 2022 (through April 26)"

# Updated / new counts
$ws.Range("R2").Value = 6
$ws.Range("B3").Value = 10
$ws.Range("V3").Value = 3
$ws.Range("F4").Value = 9
$ws.Range("N5").Value = 6
$ws.Range("V5").Value = 6
$ws.Range("AD5").Value = 2
$ws.Range("B6").Value = 9
$ws.Range("Z12").Value = 2
$ws.Range("V25").Value = 2
$ws.Range("R29").Value = 1
$ws.Range("J39").Value = 1
$ws.Range("R58").Value = 1
$ws.Range("F85").Value = 2
$ws.Range("R86").Value = 1
